# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties.
# Every player row gets the same team record: 103 wins, 58 losses, 1 tie.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1 - bold font,
# border, centered/top aligned) onto the three new header cells so they
# look consistent with the rest of the header row, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every player row (2 through 47).
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 103   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 58    # AE - Losses
    $ws.Cells.Item($r, 32).Value = 1     # AF - Ties
}
